# Fixed reversed combustion eff. calculation bug
# On the "c clinker kiln" sheet, row 7 (the Combustion / combustEff
# calculation) had its Known/Unknown quantity columns reversed:
#   before: A7=energyFuel B7=tmp   C7=fuel       D7=input
#   after:  A7=fuel       B7=input C7=energyFuel D7=tmp
# i.e. the (A,B) pair and the (C,D) pair are swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("c clinker kiln")

$knownQty   = $ws.Range("A7").Value2
$kQtyFrom   = $ws.Range("B7").Value2
$unknownQty = $ws.Range("C7").Value2
$uQtyTo     = $ws.Range("D7").Value2

$ws.Range("A7").Value = $unknownQty
$ws.Range("B7").Value = $uQtyTo
$ws.Range("C7").Value = $knownQty
$ws.Range("D7").Value = $kQtyFrom

$ws.Range("D7").Select()
